$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1017.63635
$ws.Range("I5").Value = 2078.4
$ws.Range("J5").Value = 133.66667
$ws.Range("K5").Value = 2078.4
$ws.Range("L5").Value = 133.66667
$ws.Range("M5").Value = -1963.4
$ws.Range("N5").Value = -363.66667
$ws.Range("H28").Value = 627.52
$ws.Range("I28").Value = 744.5
$ws.Range("J28").Value = 159.6
$ws.Range("K28").Value = 744.5
$ws.Range("L28").Value = 159.6
$ws.Range("M28").Value = -259.5
$ws.Range("N28").Value = -1129.6
$ws.Range("H40").Value = 6251.926
$ws.Range("I40").Value = 7933.0586
$ws.Range("J40").Value = 3394
$ws.Range("K40").Value = 7933.0586
$ws.Range("L40").Value = 3394
$ws.Range("M40").Value = -7758.0586
$ws.Range("N40").Value = -3744
$ws.Range("H117").Value = 46704.5
$ws.Range("J117").Value = 46704.5
$ws.Range("L117").Value = 46704.5
$ws.Range("N117").Value = -55882.5
$ws.Range("H130").Value = 43697.6
$ws.Range("J130").Value = 43697.6
$ws.Range("L130").Value = 43697.6
$ws.Range("N130").Value = -53737.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2916.98
$ws.Range("I32").Value = 2584.5154
$ws.Range("K32").Value = 2584.5154
$ws.Range("M32").Value = -2297.5154
$ws.Range("H122").Value = 3726.4
$ws.Range("I122").Value = 3726.4
$ws.Range("K122").Value = 11179.2
$ws.Range("M122").Value = -8729.200000000001
$ws.Range("H123").Value = 49671
$ws.Range("J123").Value = 49671
$ws.Range("L123").Value = 49671
$ws.Range("N123").Value = -59471
$ws.Range("H124").Value = 40429
$ws.Range("J124").Value = 40429
$ws.Range("L124").Value = 40429
$ws.Range("N124").Value = -50249
$ws.Range("H132").Value = 12197273
$ws.Range("I132").Value = 20834806
$ws.Range("J132").Value = 3109.2354
$ws.Range("K132").Value = 62504418
$ws.Range("L132").Value = 9327.706200000001
$ws.Range("M132").Value = -62501888
$ws.Range("N132").Value = -14387.7062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 50992
$ws.Range("J124").Value = 50992
$ws.Range("L124").Value = 50992
$ws.Range("N124").Value = -60812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 1619.4
$ws.Range("I86").Value = 1662.3077
$ws.Range("J86").Value = 1572.9166
$ws.Range("K86").Value = 1662.3077
$ws.Range("L86").Value = 1572.9166
$ws.Range("M86").Value = -539.3077000000001
$ws.Range("N86").Value = -3818.9166
$ws.Range("H89").Value = 1619.4
$ws.Range("I89").Value = 1662.3077
$ws.Range("J89").Value = 1572.9166
$ws.Range("K89").Value = 8311.538500000001
$ws.Range("L89").Value = 7864.583000000001
$ws.Range("M89").Value = -2695.538500000001
$ws.Range("N89").Value = -19096.583
$ws.Range("H100").Value = 47776
$ws.Range("J100").Value = 47776
$ws.Range("L100").Value = 47776
$ws.Range("N100").Value = -49940
$ws.Range("H116").Value = 47814.332
$ws.Range("J116").Value = 47814.332
$ws.Range("L116").Value = 47814.332
$ws.Range("N116").Value = -56992.332
$ws.Range("H132").Value = 53034.645
$ws.Range("I132").Value = 1841.7222
$ws.Range("J132").Value = 145181.9
$ws.Range("K132").Value = 5525.1666
$ws.Range("L132").Value = 435545.7
$ws.Range("M132").Value = -2995.1666
$ws.Range("N132").Value = -440605.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 22790.111
$ws.Range("I107").Value = 20580.6
$ws.Range("K107").Value = 61741.8
$ws.Range("M107").Value = -59821.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3322.182
$ws.Range("I102").Value = 3314.5
$ws.Range("J102").Value = 3342.6667
$ws.Range("K102").Value = 3314.5
$ws.Range("L102").Value = 3342.6667
$ws.Range("M102").Value = -1692.5
$ws.Range("N102").Value = -6586.6667
$ws.Range("H110").Value = 42851
$ws.Range("J110").Value = 42851
$ws.Range("L110").Value = 42851
$ws.Range("N110").Value = -51031

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H74").Value = 20216.666
$ws.Range("J74").Value = 20216.666
$ws.Range("L74").Value = 20216.666
$ws.Range("N74").Value = -22212.666
$ws.Range("H77").Value = 20216.666
$ws.Range("J77").Value = 20216.666
$ws.Range("L77").Value = 60649.99800000001
$ws.Range("N77").Value = -70633.99800000001
$ws.Range("H98").Value = 39398.4
$ws.Range("J98").Value = 39398.4
$ws.Range("L98").Value = 39398.4
$ws.Range("N98").Value = -45388.4
$ws.Range("H99").Value = 27998.666
$ws.Range("J99").Value = 29998.4
$ws.Range("L99").Value = 29998.4
$ws.Range("N99").Value = -35988.4
$ws.Range("H100").Value = 2024.5
$ws.Range("I100").Value = 1938.6428
$ws.Range("J100").Value = 2325
$ws.Range("K100").Value = 1938.6428
$ws.Range("L100").Value = 2325
$ws.Range("M100").Value = -1397.6428
$ws.Range("N100").Value = -3407
$ws.Range("H132").Value = 6627.643
$ws.Range("I132").Value = 10079.6
$ws.Range("J132").Value = 4709.8887
$ws.Range("K132").Value = 30238.8
$ws.Range("L132").Value = 14129.6661
$ws.Range("M132").Value = -27708.8
$ws.Range("N132").Value = -19189.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2908.182
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 2949
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 2949
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -4197
$ws.Range("H65").Value = 2908.182
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 2949
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 14745
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -20985
$ws.Range("H68").Value = 40997
$ws.Range("J68").Value = 40997
$ws.Range("L68").Value = 40997
$ws.Range("N68").Value = -42619
$ws.Range("H69").Value = 34999.332
$ws.Range("J69").Value = 34999.332
$ws.Range("L69").Value = 34999.332
$ws.Range("N69").Value = -36497.332
$ws.Range("H70").Value = 30104.5
$ws.Range("J70").Value = 30104.5
$ws.Range("L70").Value = 30104.5
$ws.Range("N70").Value = -30734.5
$ws.Range("H71").Value = 40997
$ws.Range("J71").Value = 40997
$ws.Range("L71").Value = 122991
$ws.Range("N71").Value = -131103
$ws.Range("H72").Value = 34999.332
$ws.Range("J72").Value = 34999.332
$ws.Range("L72").Value = 104997.996
$ws.Range("N72").Value = -112485.996
$ws.Range("H73").Value = 30104.5
$ws.Range("J73").Value = 30104.5
$ws.Range("L73").Value = 30104.5
$ws.Range("N73").Value = -32288.5
$ws.Range("H92").Value = 21766.666
$ws.Range("J92").Value = 21766.666
$ws.Range("L92").Value = 21766.666
$ws.Range("N92").Value = -26758.666
$ws.Range("H93").Value = 32466.428
$ws.Range("J93").Value = 32466.428
$ws.Range("L93").Value = 32466.428
$ws.Range("N93").Value = -37458.428
$ws.Range("H94").Value = 13557.5
$ws.Range("J94").Value = 13557.5
$ws.Range("L94").Value = 13557.5
$ws.Range("N94").Value = -15359.5
$ws.Range("H100").Value = 979.8
$ws.Range("I100").Value = 900
$ws.Range("K100").Value = 1800
$ws.Range("M100").Value = -1259
$ws.Range("H108").Value = 44208.668
$ws.Range("J108").Value = 44208.668
$ws.Range("L108").Value = 44208.668
$ws.Range("N108").Value = -51888.668
$ws.Range("H112").Value = 33851
$ws.Range("J112").Value = 33851
$ws.Range("L112").Value = 33851
$ws.Range("N112").Value = -36805
$ws.Range("H113").Value = 380.77777
$ws.Range("I113").Value = 353
$ws.Range("K113").Value = 1059
$ws.Range("M113").Value = 1111
$ws.Range("H122").Value = 1299591
$ws.Range("I122").Value = 1786438
$ws.Range("K122").Value = 5359314
$ws.Range("M122").Value = -5356864
$ws.Range("H123").Value = 43473.668
$ws.Range("J123").Value = 43473.668
$ws.Range("L123").Value = 43473.668
$ws.Range("N123").Value = -53273.668
